{"js": "// Update the payment-hours paragraph: Mon-Sat generic hours -> Mon-Fri + Sat\n// specific hours, and require \"al menos\" before \"tres d\u00edas de anticipaci\u00f3n\".\nconst oldPagos =\n  \"LOS PAGOS DEBER\u00c1N REALIZARSE DE LUNES A S\u00c1BADO, ENTRE LAS 8:30 A.M. Y LAS 5:30 P.M. \" +\n  \"PARA EFECTUAR UN PAGO EN DOMINGO, SER\u00c1 INDISPENSABLE PROGRAMAR UNA CITA CON TRES D\u00cdAS DE ANTICIPACI\u00d3N. \" +\n  \"CADA PAGO DEBER\u00c1 SER NOTIFICADO Y CONFIRMADO AL N\u00daMERO TELEF\u00d3NICO 951 189 9298.\";\nconst newPagos =\n  \"LOS PAGOS DEBER\u00c1N REALIZARSE DE LUNES A VIERNES, EN UN HORARIO DE 9:00 A.M. A 5:00 P.M., \" +\n  \"Y EN S\u00c1BADO DE 9:00 A. M. A 2:00 P. M. PARA EFECTUAR UN PAGO EN DOMINGO, SER\u00c1 INDISPENSABLE \" +\n  \"PROGRAMAR UNA CITA CON AL MENOS TRES D\u00cdAS DE ANTICIPACI\u00d3N. \" +\n  \"CADA PAGO DEBER\u00c1 SER NOTIFICADO Y CONFIRMADO AL N\u00daMERO TELEF\u00d3NICO 951 189 9298.\";\n\nconst pagosResults = context.document.body.search(oldPagos, { matchCase: true });\npagosResults.load(\"text\");\nawait context.sync();\n\nif (pagosResults.items.length > 0) {\n  pagosResults.items[0].insertText(newPagos, \"Replace\");\n  await context.sync();\n}\n\n// Update the penalty-clause paragraph: append the promisor-gender placeholder\n// clause after \"CORRESPONDAN\" (the trailing period moves to the very end).\nconst oldPena =\n  \"CUBRIR LA PENA CONVENCIONAL ESTABLECIDA POR INCUMPLIMIENTO, SIN PERJUICIO DE OTRAS ACCIONES \" +\n  \"LEGALES QUE EN SU CASO CORRESPONDAN.\";\nconst newPena =\n  \"CUBRIR LA PENA CONVENCIONAL ESTABLECIDA POR INCUMPLIMIENTO, SIN PERJUICIO DE OTRAS ACCIONES \" +\n  \"LEGALES QUE EN SU CASO CORRESPONDAN {{SEXO_7}}PROMITENTE {{SEXO_2}}\u201d.\";\n\nconst penaResults = context.document.body.search(oldPena, { matchCase: true });\npenaResults.load(\"text\");\nawait context.sync();\n\nif (penaResults.items.length > 0) {\n  penaResults.items[0].insertText(newPena, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Update the payment-hours paragraph: Mon-Sat generic hours -> Mon-Fri + Sat\n# specific hours, and require \"al menos\" before \"tres d\u00edas de anticipaci\u00f3n\".\n$d = $word.ActiveDocument\n\n$oldPagos = \"LOS PAGOS DEBER\u00c1N REALIZARSE DE LUNES A S\u00c1BADO, ENTRE LAS 8:30 A.M. Y LAS 5:30 P.M. PARA EFECTUAR UN PAGO EN DOMINGO, SER\u00c1 INDISPENSABLE PROGRAMAR UNA CITA CON TRES D\u00cdAS DE ANTICIPACI\u00d3N. CADA PAGO DEBER\u00c1 SER NOTIFICADO Y CONFIRMADO AL N\u00daMERO TELEF\u00d3NICO 951 189 9298.\"\n$newPagos = \"LOS PAGOS DEBER\u00c1N REALIZARSE DE LUNES A VIERNES, EN UN HORARIO DE 9:00 A.M. A 5:00 P.M., Y EN S\u00c1BADO DE 9:00 A. M. A 2:00 P. M. PARA EFECTUAR UN PAGO EN DOMINGO, SER\u00c1 INDISPENSABLE PROGRAMAR UNA CITA CON AL MENOS TRES D\u00cdAS DE ANTICIPACI\u00d3N. CADA PAGO DEBER\u00c1 SER NOTIFICADO Y CONFIRMADO AL N\u00daMERO TELEF\u00d3NICO 951 189 9298.\"\n\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Execute($oldPagos, $false, $true, $false, $false, $false, $true, 1, $false, $newPagos, 2) | Out-Null\n\n# Update the penalty-clause paragraph: append the promisor-gender placeholder\n# clause after \"CORRESPONDAN\" (the trailing period moves to the very end).\n$oldPena = \"CUBRIR LA PENA CONVENCIONAL ESTABLECIDA POR INCUMPLIMIENTO, SIN PERJUICIO DE OTRAS ACCIONES LEGALES QUE EN SU CASO CORRESPONDAN.\"\n$newPena = \"CUBRIR LA PENA CONVENCIONAL ESTABLECIDA POR INCUMPLIMIENTO, SIN PERJUICIO DE OTRAS ACCIONES LEGALES QUE EN SU CASO CORRESPONDAN {{SEXO_7}}PROMITENTE {{SEXO_2}}\u201d.\"\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute($oldPena, $false, $true, $false, $false, $false, $true, 1, $false, $newPena, 2) | Out-Null\n"}
